# Auto-generated edit script: updates symbol list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'269.94"
$ws.Range("E2").Value = "'3.17%"
$ws.Range("G2").Value = "'22"
$ws.Range("E3").Value = "'-1.47%"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'4.719"
$ws.Range("E4").Value = "'0.22%"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.06117"
$ws.Range("E5").Value = "'-1.51%"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'6.746"
$ws.Range("E6").Value = "'0.32%"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'0.8562"
$ws.Range("E7").Value = "'0.74%"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'0.8918"
$ws.Range("E8").Value = "'-2.17%"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'0.1432"
$ws.Range("E9").Value = "'1.83%"
$ws.Range("G9").Value = "'22"
$ws.Range("D10").Value = "'0.05042"
$ws.Range("E10").Value = "'7.90%"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.07095"
$ws.Range("E11").Value = "'0.01%"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'0.03181"
$ws.Range("E12").Value = "'0.36%"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'0.09026"
$ws.Range("E13").Value = "'-0.37%"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.001535"
$ws.Range("E14").Value = "'-0.66%"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'0.0006079"
$ws.Range("E15").Value = "'-1.10%"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.006012"
$ws.Range("E16").Value = "'-1.52%"
$ws.Range("G16").Value = "'22"
$ws.Range("E17").Value = "'-0.13%"
$ws.Range("G17").Value = "'22"
$ws.Range("D18").Value = "'3.177"
$ws.Range("E18").Value = "'0.28%"
$ws.Range("G18").Value = "'22"
$ws.Range("E19").Value = "'3.96%"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'0.3090"
$ws.Range("E20").Value = "'-0.53%"
$ws.Range("G20").Value = "'22"
$ws.Range("D21").Value = "'0.1282"
$ws.Range("E21").Value = "'-1.38%"
$ws.Range("G21").Value = "'22"
$ws.Range("D22").Value = "'3.845"
$ws.Range("E22").Value = "'-6.71%"
$ws.Range("G22").Value = "'22"
$ws.Range("D23").Value = "'0.04227"
$ws.Range("E23").Value = "'-0.01%"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("E24").Value = "'-3.09%"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.004148"
$ws.Range("E25").Value = "'0.25%"
$ws.Range("G25").Value = "'22"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("G26").Value = "'22"
$ws.Range("D27").Value = "'0.0001681"
$ws.Range("E27").Value = "'3.97%"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("G38").Value = "'22"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.03944"
$ws.Range("E40").Value = "'1.02%"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.1117"
$ws.Range("E41").Value = "'0.27%"
$ws.Range("G41").Value = "'22"
$ws.Range("D42").Value = "'0.004198"
$ws.Range("E42").Value = "'1.67%"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.002036"
$ws.Range("E43").Value = "'-6.75%"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.01327"
$ws.Range("E44").Value = "'-4.65%"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.00005139"
$ws.Range("E45").Value = "'-0.65%"
$ws.Range("G45").Value = "'22"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("G46").Value = "'22"
$ws.Range("B47").Value = "'BOLO"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.9338"
$ws.Range("E47").Value = "'457.57%"
$ws.Range("G47").Value = "'22"
$ws.Range("B48").Value = "'CoinbaseStockToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.02448"
$ws.Range("E48").Value = "'-31.83%"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("G49").Value = "'22"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("G50").Value = "'22"
$ws.Range("G51").Value = "'22"
